$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 43801
$ws.Range("A6").NumberFormat = "m/d;@"

$ws.Range("B6").Value = "黃偉倫"
$ws.Range("C6").Value = "看paper如何利用SVR技巧去把量化資料做成VIX"
$ws.Range("D6").Value = "code完成"
$ws.Range("F6").Value = "將實際資料帶入code測試"
$ws.Range("E6").Value = "初步code完成"

$ws.Range("E19").Select() | Out-Null
